$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1: 99.95 -> 0M
$t.Cell(1, 1).Range.Text = "0M"

# Row 2: 0.01 -> 0M
$t.Cell(2, 1).Range.Text = "0M"

# Row 3: 19 -> 0M
$t.Cell(3, 1).Range.Text = "0M"

# Row 4: 10 -> 40
$t.Cell(4, 1).Range.Text = "40"

# Row 6: 0.00011 -> 0.00055
$t.Cell(6, 1).Range.Text = "0.00055"

# Row 7: 0.00009 -> 0.00023
$t.Cell(7, 1).Range.Text = "0.00023"

# Row 8: 0.00003 -> 0.00007
$t.Cell(8, 1).Range.Text = "0.00007"

# Row 9: 0.00009 -> 0.00039
$t.Cell(9, 1).Range.Text = "0.00039"

# Row 10: 0.00011 -> 0.00045
$t.Cell(10, 1).Range.Text = "0.00045"

# Row 11: 0.00011 -> 0.00049
$t.Cell(11, 1).Range.Text = "0.00049"

# Row 12: 0.00095 -> 0.00932
$t.Cell(12, 1).Range.Text = "0.00932"

# Row 44: collapse multi-run tab-separated content down to "99.95"
$t.Cell(44, 1).Range.Text = "99.95"

# Row 45: collapse multi-run tab-separated content down to "0.01"
$t.Cell(45, 1).Range.Text = "0.01"

# Row 46: collapse multi-run tab-separated content down to "19"
$t.Cell(46, 1).Range.Text = "19"
